$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 457, pushing existing rows 457-543 down to 458-544.
$ws.Rows("457:457").Insert()

# Populate the newly inserted row 457 with the new weekly price record.
$ws.Cells.Item(457, 1).Value = 4
$ws.Cells.Item(457, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(457, 3).Value = "Los Lagos"
$ws.Cells.Item(457, 4).Value = 45211
$ws.Cells.Item(457, 5).Value = 10
$ws.Cells.Item(457, 6).Value = 100114014
$ws.Cells.Item(457, 7).Value = "Betarraga"
$ws.Cells.Item(457, 8).Value = "Sin especificar"
$ws.Cells.Item(457, 9).Value = "Primera"
$ws.Cells.Item(457, 10).Value = 500
$ws.Cells.Item(457, 11).Value = 1000
$ws.Cells.Item(457, 12).Value = 1000
$ws.Cells.Item(457, 13).Value = 1000
$ws.Cells.Item(457, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(457, 15).Value = "Región Metropolitana"
$ws.Cells.Item(457, 16).Value = 200
$ws.Cells.Item(457, 17).Value = 5
$ws.Cells.Item(457, 18).Value = "Hortaliza"
